# Add a new "Save" column (H) to the sheet, mirroring the style of the
# other header cells and filling in the per-row save indicator values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 -------------------------------------------------
# Copy the formatting from the neighboring header cell (G1, style s="1":
# bold font, thin border, centered/top aligned) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# --- Data values H2:H58 ----------------------------------------------
$saveValues = @(1,0,0,0,1,1,1,1,0,0,0,1,1,0,0,0,0,0,0,1,0,1,0,0,1,0,0,1,0,0,0,0,0,0,0,0,1,1,0,0,0,1,1,0,0,0,1,1,1,1,0,1,1,0,1,1,1)

$startRow = 2
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}

Write-Host "Save column added"
